$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

# Column A holds dates stored as plain text in this sheet (e.g. "2025-05-02"),
# so force text formatting before assigning the value to avoid Excel's
# automatic date-serial conversion, then restore the default "Normal" style
# so the new cell doesn't pick up a stray number format.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-05-05"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item($row, 3).Value = "NA"
$ws.Cells.Item($row, 4).Value = 1
